$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 613.5
$ws.Range("J5").Value = 1067
$ws.Range("L5").Value = 1067
$ws.Range("N5").Value = -1297
$ws.Range("H7").Value = 14321.667
$ws.Range("I7").Value = 2965
$ws.Range("J7").Value = 20000
$ws.Range("K7").Value = 2965
$ws.Range("L7").Value = 20000
$ws.Range("M7").Value = -2853
$ws.Range("N7").Value = -20224
$ws.Range("H14").Value = 14321.667
$ws.Range("I14").Value = 2965
$ws.Range("J14").Value = 20000
$ws.Range("K14").Value = 2965
$ws.Range("L14").Value = 20000
$ws.Range("M14").Value = -2774
$ws.Range("N14").Value = -20382
$ws.Range("H17").Value = 2158
$ws.Range("J17").Value = 2158
$ws.Range("L17").Value = 6474
$ws.Range("N17").Value = -6810
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H70").Value = 1462.1428
$ws.Range("I70").Value = 1422.5
$ws.Range("K70").Value = 4267.5
$ws.Range("M70").Value = -3997.5
$ws.Range("H73").Value = 1462.1428
$ws.Range("I73").Value = 1422.5
$ws.Range("K73").Value = 4267.5
$ws.Range("M73").Value = -3331.5
$ws.Range("H80").Value = 556.1667
$ws.Range("I80").Value = 635.75
$ws.Range("J80").Value = 397
$ws.Range("K80").Value = 1907.25
$ws.Range("L80").Value = 1191
$ws.Range("M80").Value = -909.25
$ws.Range("N80").Value = -3187
$ws.Range("H81").Value = 70000
$ws.Range("J81").Value = 70000
$ws.Range("L81").Value = 70000
$ws.Range("N81").Value = -71996
$ws.Range("H83").Value = 556.1667
$ws.Range("I83").Value = 635.75
$ws.Range("J83").Value = 397
$ws.Range("K83").Value = 5721.75
$ws.Range("L83").Value = 3573
$ws.Range("M83").Value = -729.75
$ws.Range("N83").Value = -13557
$ws.Range("H84").Value = 70000
$ws.Range("J84").Value = 70000
$ws.Range("L84").Value = 210000
$ws.Range("N84").Value = -219984
$ws.Range("H94").Value = 979
$ws.Range("I94").Value = 979
$ws.Range("K94").Value = 979
$ws.Range("M94").Value = -528
$ws.Range("H124").Value = 150000
$ws.Range("J124").Value = 150000
$ws.Range("L124").Value = 150000
$ws.Range("N124").Value = -159820
$ws.Range("H132").Value = 1093.3846
$ws.Range("I132").Value = 953.7143
$ws.Range("J132").Value = 1680
$ws.Range("K132").Value = 2861.1429
$ws.Range("L132").Value = 5040
$ws.Range("M132").Value = -331.1428999999998
$ws.Range("N132").Value = -10100
$ws.Range("H135").Value = 897.88
$ws.Range("I135").Value = 802.087
$ws.Range("K135").Value = 7218.782999999999
$ws.Range("M135").Value = -4683.782999999999
$ws.Range("H137").Value = 3300.6667
$ws.Range("I137").Value = 2301.5
$ws.Range("J137").Value = 6498
$ws.Range("K137").Value = 6904.5
$ws.Range("L137").Value = 19494
$ws.Range("M137").Value = -4354.5
$ws.Range("N137").Value = -24594

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 97.75
$ws.Range("I5").Value = 97.75
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 97.75
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 14.25
$ws.Range("N5").ClearContents()
$ws.Range("H45").Value = 3118.1333
$ws.Range("J45").Value = 4332.6665
$ws.Range("L45").Value = 4332.6665
$ws.Range("N45").Value = -5086.6665
$ws.Range("H97").Value = 1417.8889
$ws.Range("I97").Value = 330.66666
$ws.Range("K97").Value = 330.66666
$ws.Range("M97").Value = 165.33334
$ws.Range("H110").Value = 1600
$ws.Range("I110").Value = 1600
$ws.Range("K110").Value = 1600
$ws.Range("M110").Value = 445
$ws.Range("H132").Value = 3017.5293
$ws.Range("I132").Value = 2530.7693
$ws.Range("J132").Value = 4599.5
$ws.Range("K132").Value = 7592.3079
$ws.Range("L132").Value = 13798.5
$ws.Range("M132").Value = -5062.3079
$ws.Range("N132").Value = -18858.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 203
$ws.Range("I7").Value = 203
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 203
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -90
$ws.Range("N7").ClearContents()
$ws.Range("H20").Value = 5408.0835
$ws.Range("I20").Value = 5255.4443
$ws.Range("J20").Value = 5866
$ws.Range("K20").Value = 5255.4443
$ws.Range("L20").Value = 5866
$ws.Range("M20").Value = -5008.4443
$ws.Range("N20").Value = -6360
$ws.Range("H134").Value = 4607.5
$ws.Range("I134").Value = 4607.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 13822.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -11287.5
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 682.6667
$ws.Range("J22").Value = 659.2
$ws.Range("L22").Value = 659.2
$ws.Range("N22").Value = -1359.2
$ws.Range("H31").Value = 2458.375
$ws.Range("I31").Value = 2540.9092
$ws.Range("J31").Value = 2276.8
$ws.Range("K31").Value = 2540.9092
$ws.Range("L31").Value = 2276.8
$ws.Range("M31").Value = -2245.9092
$ws.Range("N31").Value = -2866.8
$ws.Range("H34").Value = 2458.375
$ws.Range("I34").Value = 2540.9092
$ws.Range("J34").Value = 2276.8
$ws.Range("K34").Value = 2540.9092
$ws.Range("L34").Value = 2276.8
$ws.Range("M34").Value = -2338.9092
$ws.Range("N34").Value = -2680.8
$ws.Range("H58").Value = 3006
$ws.Range("I58").Value = 2998
$ws.Range("K58").Value = 2998
$ws.Range("M58").Value = -2795
$ws.Range("H105").Value = 6050
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("H107").Value = 1132.6666
$ws.Range("I107").Value = 390.83334
$ws.Range("K107").Value = 390.83334
$ws.Range("M107").Value = 1529.16666
$ws.Range("H132").Value = 3022.25
$ws.Range("I132").Value = 2881.6428
$ws.Range("J132").Value = 4006.5
$ws.Range("K132").Value = 8644.928400000001
$ws.Range("L132").Value = 12019.5
$ws.Range("M132").Value = -6114.928400000001
$ws.Range("N132").Value = -17079.5
$ws.Range("H134").Value = 2939.3333
$ws.Range("I134").Value = 3015.6365
$ws.Range("J134").Value = 2100
$ws.Range("K134").Value = 9046.9095
$ws.Range("L134").Value = 6300
$ws.Range("M134").Value = -6511.9095
$ws.Range("N134").Value = -11370
$ws.Range("H136").Value = 3006
$ws.Range("I136").Value = 2998
$ws.Range("K136").Value = 8994
$ws.Range("M136").Value = -6444

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 229.66667
$ws.Range("J23").Value = 349
$ws.Range("L23").Value = 1047
$ws.Range("N23").Value = -1517
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()
$ws.Range("H98").Value = 14348.5
$ws.Range("J98").Value = 20748.5
$ws.Range("L98").Value = 62245.5
$ws.Range("N98").Value = -65241.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 89.40000000000001
$ws.Range("I2").Value = 74.375
$ws.Range("K2").Value = 74.375
$ws.Range("M2").Value = 38.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2167500
$ws.Range("I2").Value = 2167500
$ws.Range("K2").Value = 2167500
$ws.Range("M2").Value = -2167388
$ws.Range("H7").Value = 3098
$ws.Range("I7").Value = 3098
$ws.Range("K7").Value = 3098
$ws.Range("M7").Value = -2986
$ws.Range("H16").Value = 777
$ws.Range("I16").Value = 777
$ws.Range("K16").Value = 777
$ws.Range("M16").Value = -607
$ws.Range("H93").Value = 2272
$ws.Range("I93").Value = 2272
$ws.Range("K93").Value = 2272
$ws.Range("M93").Value = -1024
$ws.Range("H126").Value = 3098
$ws.Range("I126").Value = 3098
$ws.Range("K126").Value = 9294
$ws.Range("M126").Value = -6824
$ws.Range("H132").Value = 3431.6
$ws.Range("I132").Value = 2664
$ws.Range("J132").Value = 3623.5
$ws.Range("K132").Value = 7992
$ws.Range("L132").Value = 10870.5
$ws.Range("M132").Value = -5462
$ws.Range("N132").Value = -15930.5
$ws.Range("H136").Value = 3017.5925
$ws.Range("I136").Value = 2858.6
$ws.Range("J136").Value = 5005
$ws.Range("K136").Value = 8575.799999999999
$ws.Range("L136").Value = 15015
$ws.Range("M136").Value = -6025.799999999999
$ws.Range("N136").Value = -20115

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 25000
$ws.Range("J2").Value = 25000
$ws.Range("L2").Value = 25000
$ws.Range("N2").Value = -25224
$ws.Range("H132").Value = 2351.1462
$ws.Range("I132").Value = 1553.7587
$ws.Range("K132").Value = 4661.2761
$ws.Range("M132").Value = -2131.2761
$ws.Range("H136").Value = 999.2222
$ws.Range("I136").Value = 922.46155
$ws.Range("J136").Value = 2995
$ws.Range("K136").Value = 2767.38465
$ws.Range("L136").Value = 8985
$ws.Range("M136").Value = -217.38465
$ws.Range("N136").Value = -14085
